# figure_02: panel tags to uppercase, darken the "bright landscapes" label

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "bright landscapes" label -> darker gold (F7C267 -> F4A820)
$brightLandscapes = $s.Shapes.Item(16)
$brightLandscapes.TextFrame.TextRange.Font.Color.RGB = 2140404  # RGB(0xF4,0xA8,0x20)

# Panel tags (a)/(b)/(c) -> uppercase (A)/(B)/(C)
$panelA = $s.Shapes.Item(34)
$panelA.TextFrame.TextRange.Text = "(A)"

$panelB = $s.Shapes.Item(40)
$panelB.TextFrame.TextRange.Text = "(B)"

$panelC = $s.Shapes.Item(41)
$panelC.TextFrame.TextRange.Text = "(C)"
